{"js": "// Insert a new bulleted list item \"Sollte keine Anrede passen, so wird das\n// Geschlecht \u201edivers\" bzw. \u201eneutral\" eingetragen\" right after the existing\n// list item \"Anrede definiert eindeutig Geschlecht\", matching its list\n// formatting (style \"Listenabsatz\", bullet list numId 1).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the anchor paragraph by its text.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Anrede definiert eindeutig Geschlecht\") {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not locate paragraph \"Anrede definiert eindeutig Geschlecht\"');\n}\n\n// insertParagraph(\"After\") splits off a new paragraph that inherits the\n// anchor's paragraph formatting (style + list numbering), then we fill in\n// its text.\nanchor.insertParagraph(\n  \"Sollte keine Anrede passen, so wird das Geschlecht \u201edivers\u201c bzw. \u201eneutral\u201c eingetragen\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Insert a new bulleted list item \"Sollte keine Anrede passen, so wird das\n# Geschlecht \u201edivers\" bzw. \u201eneutral\" eingetragen\" right after the existing\n# list item \"Anrede definiert eindeutig Geschlecht\", matching its list\n# formatting (style \"Listenabsatz\", bullet list numId 1).\n\n$d = $word.ActiveDocument\n\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`n\") -eq \"Anrede definiert eindeutig Geschlecht\") {\n        $anchor = $p\n        break\n    }\n}\n\nif ($null -eq $anchor) {\n    throw 'Could not locate paragraph \"Anrede definiert eindeutig Geschlecht\"'\n}\n\n# InsertParagraphAfter splits a new paragraph in after the anchor, inheriting\n# its paragraph formatting (style + list numbering). It returns nothing (a\n# VBA Sub), so re-fetch the freshly created (now-next) paragraph afterwards\n# and set its text.\n$anchor.Range.InsertParagraphAfter() | Out-Null\n$newParagraph = $anchor.Next()\n$newParagraph.Range.Text = \"Sollte keine Anrede passen, so wird das Geschlecht \u201edivers\u201c bzw. \u201eneutral\u201c eingetragen\"\n"}
